$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text in the Price/Volume columns is preserved as literal
# text (matching the source data export) rather than auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "91.255.59"
$ws.Range("E2").Value = "  +4.47%  "
$ws.Range("D3").Value = "3.207.40"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "221.12"
$ws.Range("E5").Value = "  +7.07%  "
$ws.Range("D6").Value = "641.72"
$ws.Range("E6").Value = "  +6.18%  "
$ws.Range("D7").Value = "0.404"
$ws.Range("E7").Value = "  +6.94%  "
$ws.Range("D8").Value = "0.713"
$ws.Range("E8").Value = "  +7.96%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "3.205.02"
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("D11").Value = "0.578"
$ws.Range("E11").Value = "  +9.20%  "
$ws.Range("D12").Value = "0.182"
$ws.Range("E12").Value = "  +3.83%  "
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  +8.80%  "
$ws.Range("D14").Value = "5.44"
$ws.Range("E14").Value = "  +4.85%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "90.926.90"
$ws.Range("E15").Value = "  +4.32%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "33.47"
$ws.Range("E16").Value = "  +5.40%  "
$ws.Range("D17").Value = "3.799.84"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "3.218.69"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("D19").Value = "0.0000232"
$ws.Range("E19").Value = "  +80.04%  "
$ws.Range("E20").Value = "  +8.96%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "442.79"
$ws.Range("E21").Value = "  +7.89%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "13.48"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("E23").Value = "  +3.54%  "
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").Value = "5.37"
$ws.Range("E25").Value = "  +5.24%  "
$ws.Range("D26").Value = "11.94"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").Value = "81.61"
$ws.Range("E27").Value = "  +12.16%  "
$ws.Range("D28").Value = "3.378.99"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "0.160"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").Value = "4.23"
$ws.Range("E32").Value = "  +42.50%  "
$ws.Range("D33").Value = "8.48"
$ws.Range("E33").Value = "  +4.38%  "
$ws.Range("D34").Value = "537.24"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "7.13"
$ws.Range("E35").Value = "  +8.00%  "
$ws.Range("E36").Value = "  +4.74%  "
$ws.Range("D37").Value = "1.30"
$ws.Range("E37").Value = "  +1.94%  "
$ws.Range("D38").Value = "22.62"
$ws.Range("E38").Value = "  +4.74%  "
$ws.Range("D39").Value = "22.40"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("E42").Value = "  +3.73%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("D45").Value = "147.54"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Value = "44.86"
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("D47").Value = "174.33"
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("D48").Value = "0.127"
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("D49").Value = "0.757"
$ws.Range("E49").Value = "  +10.25%  "
$ws.Range("D50").Value = "25.40"
$ws.Range("E50").Value = "  +8.55%  "
$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").Value = "1.24"
$ws.Range("E51").Value = "  +2.91%  "

# Restore default cell style on the Price/Volume columns (clears the temporary
# text-format override applied above) without touching the stored text values.
$ws.Range("D2:E51").Style = "Normal"

Write-Host "Updated cryptos list"
